$d = $word.ActiveDocument

$map = @{
    "Имя регистра" = "Регистр";
    "Имя поля"     = "Поле";
    "Reset"        = "Значение";
    "Имя enum"     = "Enum";
}

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $t = $r.Text
    $trimmed = $t.TrimEnd("`r", "`n", "`a", "`v", "`f")
    if ($map.ContainsKey($trimmed)) {
        $r.Text = $map[$trimmed]
    }
}
